# The sheet originally had two rows:
#   A1 = 0            (bold font, thin border -> style index 1, a "header" cell)
#   A2 = <big string>  (default style) holding the questions payload as a
#                       Python-literal dict/list dump.
#
# The edit removes the standalone header cell entirely (so the payload moves
# up to A1 with the plain/default style) and reformats the payload text from
# a single-quoted Python literal into pretty-printed JSON (double quotes,
# 4-space indent, one entry per line).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop row 1 (the bold/bordered "0" header). This shifts the payload that
# used to live in A2 up into A1, picking up the sheet's plain default style.
$ws.Rows(1).Delete()

# Overwrite the cell text with the pretty-printed JSON rendering of the same
# data (the list of question dicts, now serialized as JSON with indent=4).
$ws.Range("A1").Value = 'questions = [
    {
        "title": "The table below presents instructions paired with corresponding images. Each instruction should correctly match the image to which it is paired. However, one of the pairings is incorrect, as the instructions do not match the image presented. Which pairing is incorrect?",
        "ques_type": 2,
        "options": [
            "A",
            "B",
            "C",
            "D"
        ],
        "score": "A"
    },
    {
        "title": "As a warehouse picker, you have gathered all items according to your pickup list and are labeling them for shipment. During this process, you notice that one package has torn corners.What should you do next?",
        "ques_type": 2,
        "options": [
            "Repack the item yourself, using available packing material.",
            "Give the package to the packing team for repacking.",
            "Assess whether the package requires repacking.",
            "Secure the torn corners with packaging tape."
        ],
        "score": "Assess whether the package requires repacking."
    },
    {
        "title": "You manage storage in an automobile warehouse. Four sections are in the warehouse, including the following:An open section for large-sized items that cannot be put into cartons, such as machineryA rack section for medium-sized items that can be put in cartons, such as automobile spare partsA bin section for small items, such as screwsA cool, dry warehouse section for temperature-sensitive chemicals, such as paint and resinWhere should you store the paint thinner?",
        "ques_type": 2,
        "options": [
            "Open section",
            "Rack section",
            "Bin section",
            "Cool, dry section"
        ],
        "score": "Cool, dry section"
    }
]'

# Setting a value containing embedded line breaks makes the engine pin an
# explicit (custom) row height to fit all the wrapped lines; AutoFit the row
# back so no custom height is recorded, matching a plain, unstyled row.
$ws.Rows(1).AutoFit()
